# Add the "I0" and "IF" columns (I and J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - reuse the same bold/bordered style as the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-64.
$i0 = @(5,8,6,8,6,1,6,4,7,5,9,8,6,5,8,8,8,8,5,6,4,5,7,7,7,1,7,6,1,6,9,6,5,6,7,5,7,6,6,1,7,5,8,6,5,4,7,6,7,8,9,4,8,7,8,5,7,9,2,8,4,5,3)
$if = @(6,8,6,8,6,4,7,5,8,7,9,8,7,7,8,8,9,8,6,7,6,6,7,7,7,4,7,7,3,7,9,8,5,7,7,5,8,6,6,3,7,6,8,7,5,6,7,6,7,8,9,5,9,8,8,6,7,9,5,8,4,6,3)

for ($idx = 0; $idx -lt $i0.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0[$idx]
    $ws.Cells.Item($row, 10).Value = $if[$idx]
}
